# "addind preparer to sheet"
# Update the libraryPreparer column (B) to the actual preparer's initials,
# and rename the purpose column (E) value from the placeholder
# "Retrofitted_4019" to "fullRNASEQ" for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 5).Value = "fullRNASEQ"
    $ws.Cells.Item($r, 2).Value = "H.BROWN"
}

$ws.Range("B3:B19").Select()
